# Update gh-pages to output generated at 456a3b4
# Applies numeric / content refreshes to the "展览" (sheet1), "演出" (sheet2)
# and "全部类型" (sheet4) worksheets of the 合肥-漫展信息 workbook.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G2").Value = "不可售"
$ws1.Range("F3").Value = 527
$ws1.Range("F4").Value = 40
$ws1.Range("F5").Value = 25
$ws1.Range("F8").Value = 29
$ws1.Range("F9").Value = 283
$ws1.Range("F10").Value = 2988
$ws1.Range("F11").Value = 30

# ---- Sheet "演出" ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 45

# ---- Sheet "全部类型" ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value = "不可售"
$ws4.Range("F4").Value = 527
$ws4.Range("F5").Value = 40
$ws4.Range("F6").Value = 25
$ws4.Range("F9").Value = 29
$ws4.Range("F10").Value = 283
$ws4.Range("F11").Value = 2988
$ws4.Range("F12").Value = 30
$ws4.Range("F13").Value = 45
